$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.517.69"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.631.89"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'112.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'324.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.544"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").Value = "'39.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").Value = "'19.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.25%  "
$ws.Range("D12").Value = "'0.0812"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "'7.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "3.041.23"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "2.622.43"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'0.849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("D18").Value = "49.421.80"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "0.0₃0947"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("D23").Value = "'269.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.54%  "
$ws.Range("D24").Value = "'68.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("D25").Value = "'2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'10.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("D29").Value = "'2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("D31").Value = "'34.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.31%  "
$ws.Range("D32").Value = "'49.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'5.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'18.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").Value = "'4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'128.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").Value = "'22.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'0.0326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.11%  "
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("D45").Value = "2.058.74"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  -4.89%  "
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("E48").Value = "  -5.41%  "
$ws.Range("D49").Value = "'8.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").Value = "'59.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("E51").Value = "  -4.29%  "
